$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("B2").Value = 19.3164392298743
$ws.Range("C2").Value = 0.02882075676116136

$ws.Range("B3").Value = 19.91364414806782
$ws.Range("C3").Value = 0.04773430819931494

$ws.Range("B4").Value = 19.24633245902404
$ws.Range("C4").Value = 0.03836751660295683
